$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Remove the now-obsolete rows (old "Check comparision/ResultType" header
# row and the three "ResultNumN" rows) working bottom-to-top so the row
# indices of rows above the deletion point stay stable. ---
$ws.Rows("12:12").Delete()
$ws.Rows("10:10").Delete()
$ws.Rows("8:8").Delete()
$ws.Rows("6:6").Delete()

# --- Header row tweak ---
$ws.Range("B4").Value = "Parameter Name"

# --- Row 6: merged "Check comparison results" label + DiffCount values ---
$ws.Range("A6").Value = "Check comparison results"
$ws.Range("B6").Value = "DiffCount"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3

# --- Row 7: ResultText0 ---
$ws.Range("B7").Value = "ResultText0"
$ws.Range("C7").Value = "A1   A1  NewA1"

# --- Row 8: ResultText1 ---
$ws.Range("B8").Value = "ResultText1"
$ws.Range("C8").Value = "B6   B6  NewB6"

# --- Row 9: ResultText2 ---
$ws.Range("B9").Value = "ResultText2"
$ws.Range("C9").Value = "C11  C11 NewC11"

# --- Wrap text on the merged label column (matches new cellXfs that add
# wrapText="1" to the fillId=4 border styles used by column A). ---
$ws.Range("A6:A9").WrapText = $true

# --- Column widths to roughly match the new best-fit layout ---
$ws.Columns.Item(1).ColumnWidth = 9.3
$ws.Columns.Item(2).ColumnWidth = 16.3
$ws.Columns.Item(3).ColumnWidth = 16
$ws.Columns.Item(4).ColumnWidth = 11.3
$ws.Columns.Item(5).ColumnWidth = 14.4
$ws.Columns.Item(6).ColumnWidth = 10

# --- Selection state matches the authored workbook ---
$ws.Range("L7").Select()
